$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 3.045497666666666
$ws.Range("H2").Value = 9.136493
$ws.Range("I2").Value = 0.06184575966423571
$ws.Range("J2").Value = 0.06184575966423572
$ws.Range("M2").Value = 1.443038
$ws.Range("N2").Value = 4.329114
$ws.Range("O2").Value = 0.0289666880885598
$ws.Range("P2").Value = 0.0289666880885598
$ws.Range("Q2").Value = 4.394768861911333
$ws.Range("R2").Value = 39.552919757202
$ws.Range("S2").Value = 0.001791466829793949
$ws.Range("T2").Value = 0.001791466829793949
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 3.045497666666666
$ws.Range("H3").Value = 9.136493
$ws.Range("I3").Value = 0.06184575966423571
$ws.Range("J3").Value = 0.06184575966423572
$ws.Range("N3").Value = 87.61054300000001
$ws.Range("O3").Value = 0.5862140087672342
$ws.Range("P3").Value = 0.5862140087672342
$ws.Range("Q3").Value = 88.93923476063321
$ws.Range("R3").Value = 800.453112845699
$ws.Range("S3").Value = 0.03625485069802653
$ws.Range("T3").Value = 0.03625485069802653
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 3.045497666666666
$ws.Range("H4").Value = 9.136493
$ws.Range("I4").Value = 0.06184575966423571
$ws.Range("J4").Value = 0.06184575966423572
$ws.Range("M4").Value = 19.170603
$ws.Range("N4").Value = 57.511809
$ws.Range("O4").Value = 0.384819303144206
$ws.Range("P4").Value = 0.384819303144206
$ws.Range("Q4").Value = 58.38402670509299
$ws.Range("R4").Value = 525.4562403458369
$ws.Range("S4").Value = 0.02379944213641523
$ws.Range("T4").Value = 0.02379944213641523
$ws.Range("I5").Value = 0.6352626115862781
$ws.Range("J5").Value = 0.6352626115862781
$ws.Range("M5").Value = 1.443038
$ws.Range("N5").Value = 4.329114
$ws.Range("O5").Value = 0.0289666880885598
$ws.Range("P5").Value = 0.0289666880885598
$ws.Range("Q5").Value = 45.14185547550667
$ws.Range("R5").Value = 406.27669927956
$ws.Range("S5").Value = 0.01840145392414363
$ws.Range("T5").Value = 0.01840145392414363
$ws.Range("I6").Value = 0.6352626115862781
$ws.Range("J6").Value = 0.6352626115862781
$ws.Range("N6").Value = 87.61054300000001
$ws.Range("O6").Value = 0.5862140087672342
$ws.Range("P6").Value = 0.5862140087672342
$ws.Range("Q6").Value = 913.5593265126913
$ws.Range("R6").Value = 8222.033938614222
$ws.Range("S6").Value = 0.3723998421579345
$ws.Range("T6").Value = 0.3723998421579345
$ws.Range("I7").Value = 0.6352626115862781
$ws.Range("J7").Value = 0.6352626115862781
$ws.Range("M7").Value = 19.170603
$ws.Range("N7").Value = 57.511809
$ws.Range("O7").Value = 0.384819303144206
$ws.Range("P7").Value = 0.384819303144206
$ws.Range("Q7").Value = 599.7046439555401
$ws.Range("R7").Value = 5397.34179559986
$ws.Range("S7").Value = 0.2444613155041999
$ws.Range("T7").Value = 0.2444613155041999
$ws.Range("G8").Value = 14.91542433333333
$ws.Range("H8").Value = 44.746273
$ws.Range("I8").Value = 0.3028916287494862
$ws.Range("J8").Value = 0.3028916287494862
$ws.Range("M8").Value = 1.443038
$ws.Range("N8").Value = 4.329114
$ws.Range("O8").Value = 0.0289666880885598
$ws.Range("P8").Value = 0.0289666880885598
$ws.Range("Q8").Value = 21.52352409912466
$ws.Range("R8").Value = 193.711716892122
$ws.Range("S8").Value = 0.008773767334622219
$ws.Range("T8").Value = 0.008773767334622219
$ws.Range("G9").Value = 14.91542433333333
$ws.Range("H9").Value = 44.746273
$ws.Range("I9").Value = 0.3028916287494862
$ws.Range("J9").Value = 0.3028916287494862
$ws.Range("N9").Value = 87.61054300000001
$ws.Range("O9").Value = 0.5862140087672342
$ws.Range("P9").Value = 0.5862140087672342
$ws.Range("Q9").Value = 435.5828083062488
$ws.Range("R9").Value = 3920.24527475624
$ws.Range("S9").Value = 0.1775593159112732
$ws.Range("T9").Value = 0.1775593159112732
$ws.Range("G10").Value = 14.91542433333333
$ws.Range("H10").Value = 44.746273
$ws.Range("I10").Value = 0.3028916287494862
$ws.Range("J10").Value = 0.3028916287494862
$ws.Range("M10").Value = 19.170603
$ws.Range("N10").Value = 57.511809
$ws.Range("O10").Value = 0.384819303144206
$ws.Range("P10").Value = 0.384819303144206
$ws.Range("Q10").Value = 285.937678470873
$ws.Range("R10").Value = 2573.439106237857
$ws.Range("S10").Value = 0.1165585455035908
$ws.Range("T10").Value = 0.1165585455035908
